$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for existing rows (2-10), columns B:G
$data = @{
    2  = @(0.02846637460976744, 0.9430602768606053, 2.648581526649609, 1.627446320666095, 1.643388759852491, 51)
    3  = @(-0.1943718813857691, 0.7420322100163707, 1.301007984511005, 1.140617369897112, 1.13534475400492, 50)
    4  = @(-0.05603653138027285, 0.7181205273575723, 1.474592598453738, 1.214328044003653, 1.22560506124038, 49)
    5  = @(-0.151782044080011, 0.7437703319808541, 1.325335475341544, 1.151232155276052, 1.153258924895003, 48)
    6  = @(-0.01202353621938327, 0.9685194911816654, 2.135832188922473, 1.461448661062876, 1.47719857003534, 47)
    7  = @(-0.1554833155668876, 0.8002215085028028, 1.375155290479728, 1.172670154169419, 1.175160419656847, 46)
    8  = @(0.04659670587843785, 0.8050767692874348, 1.154340410675184, 1.074402350460564, 1.085520529677319, 45)
    9  = @(-0.1084043684214821, 0.735808405588202, 1.20283139358349, 1.096736702031755, 1.1039834322736, 44)
    10 = @(0.1182394667038753, 1.047860666370489, 2.217418452168737, 1.489099879849816, 1.501965607916916, 43)
    11 = @(-0.08993503749128429, 0.8551623832460461, 1.503988522260922, 1.226372097799409, 1.237895626577417, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i   # B=2 ... G=7
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# New row 11 label in column A, mirroring formatting of existing A column entries (e.g. A10 "Q8")
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").Value = "Q9"
